$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 3409
$ws1.Range("F9").Value = 4084
$ws1.Range("F11").Value = 1021

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 3409
$ws4.Range("F10").Value = 4084
$ws4.Range("F12").Value = 1021
